$d = $word.ActiveDocument

$ptXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1. Introdução aos sistemas de manufatura;</w:t><w:br/><w:t xml:space="preserve">    2. Indústria 4.0;</w:t><w:br/><w:t xml:space="preserve">    3. Inteligencia artifical (ai), internet das coisas (IoT) e sistemas ciberfísicos;</w:t><w:br/><w:t xml:space="preserve">    4. Sistemas de controle industrial; conceito de PLM e integração com os sistemas de gestão;</w:t><w:br/><w:t xml:space="preserve">    5. Monitoramento e supervisão de processos de produção. Sistemas de controle da produção, manufatura sustentável;</w:t><w:br/><w:t xml:space="preserve">    6. Componentes de hardware para automação de processos: controle numérico, programação CNC, controle discreto utilizando controladores lógico programáveis e sistemas on-chip;</w:t><w:br/><w:t xml:space="preserve">    7. Robótica industrial – programação de robôs e robôs colaborativos;</w:t><w:br/><w:t xml:space="preserve">    8. Sistemas de transporte de materiais e sistemas de armazenamento;</w:t><w:br/><w:t xml:space="preserve">    9. Identificação automática e captura de dados – tecnologias de inspeção.</w:t></w:r></w:p>'
$enXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>1. Introduction to manufacturing systems;</w:t><w:br/><w:t>2. 4.0 industry;</w:t><w:br/><w:t>3. Artificial intelligence (AI), internet of things (IoT), and cyber-physical systems;</w:t><w:br/><w:t>4. Industrial control systems; PLM concept and integration with management systems;</w:t><w:br/><w:t>5. Production process monitoring and supervision. Production control systems, sustainable manufacturing;</w:t><w:br/><w:t>6. Hardware components for process automation: numerical control, CNC programming, discrete control using programmable logic controllers and on-chip systems;</w:t><w:br/><w:t>7. Industrial robotics – programming of robots and collaborative robots;</w:t><w:br/><w:t>8. Material transportation systems and storage systems;</w:t><w:br/><w:t>9. Automatic identification and data capture – inspection technologies.</w:t></w:r></w:p>'
$biblioXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>GROOVER, M.P. Automação Industrial e Sistemas de Manufatura, 561p., 3a Edição - São Paulo, Pearson Prentice Hall, 2011.</w:t><w:br/><w:t>RICHARD, L.S., ERNEST, L.H., Handbook of Industrial Automation, Marcel Dekker,Inc. NewYork, 2000.</w:t><w:br/><w:t>ADALBERTO FILHO ET. AL, Automação &amp; Sociedade: Quarta revolução Industrial, um olhar para o Brasil, 1a Edição, Brasport Livros e Multimídia Limitada.</w:t></w:r></w:p>'

foreach ($p in $d.Paragraphs) {
  $t = $p.Range.Text
  if ($t.StartsWith('1. Introdução aos sistemas de manufatura;')) {
    [void]$p.Range.InsertXML($ptXml)
  }
  elseif ($t.StartsWith('1. Introduction to manufacturing systems;')) {
    [void]$p.Range.InsertXML($enXml)
  }
  elseif ($t.StartsWith('GROOVER, M.P.')) {
    [void]$p.Range.InsertXML($biblioXml)
  }
}
